# Updated fitting parameters.
# "Parameters" sheet, column K ("h_p_star"), row 2: 0.2733 -> 0.28125

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
$ws.Range("K2").Value = 0.28125
